$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values per diff
$ws.Range("D3").Value = 0.0008441589660588925

$ws.Range("B4").Value = 100
$ws.Range("C4").Value = 0.7901649589302215
$ws.Range("D4").Value = 0.02149806927100434
$ws.Range("E4").Value = 0.8581399114944516

$ws.Range("D5").Value = 0.008493266987447642

# Add new row 6 for DWA (copy formatting from A5, which carries the bold/border/center style)
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = "DWA"

$ws.Range("B6").Value = 100
$ws.Range("C6").Value = 0.9364989339188737
$ws.Range("D6").Value = 0.02858297303144933
$ws.Range("E6").Value = 0.9498956054774071
